$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44320, 710, 2938, 550, 5601, 201, 1927, 1000, 28500, 0, 0, 0, 3638, 265, 1200, 42605),
    @(44321, 250, 2703, 870, 5525, 131, 1928, 0, 28500, 0, 0, 0, 3638, 158, 1213, 42294)
)

$startRow = 87
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
}
